# Update code tinh luong % format cac bang
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 2: "Đơn phụ phẫu 1" -- collapse the detailed columns (J..AA) into a
# much smaller table, and repurpose/rename columns G, H, I.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New header row (columns G, H, I renamed; J..AA removed)
$ws2.Range("G1").Value = "Tên dịch vụ"
$ws2.Range("H1").Value = "Phụ phẫu 1"
$ws2.Range("I1").Value = "Công phụ phẫu 1"

# Data rows 2-6: G = old "Tên dịch vụ" (H) value, H = old "Phụ phẫu 1" (T) value,
# I = old "Công phụ phẫu 1" (V) value
$ws2.Range("G2").Value = "Tiềm cằm"
$ws2.Range("H2").Value = "Đào Vương Anh"
$ws2.Range("I2").Value = 0

$ws2.Range("G3").Value = "Cắt mí"
$ws2.Range("H3").Value = "Đào Vương Anh"
$ws2.Range("I3").Value = 50000

$ws2.Range("G4").Value = "Tiêm môi"
$ws2.Range("H4").Value = "Đào Vương Anh"
$ws2.Range("I4").Value = 0

$ws2.Range("G5").Value = "Cắt mí"
$ws2.Range("H5").Value = "Đào Vương Anh"
$ws2.Range("I5").Value = 50000

$ws2.Range("G6").Value = "Thu cánh mũi"
$ws2.Range("H6").Value = "Đào Vương Anh"
$ws2.Range("I6").Value = 100000

# Totals row 7
$ws2.Range("I7").Value = 200000

# Clear out everything from column J through AA for rows 1-7 (old detail
# columns that no longer exist in the new, narrower table).
$ws2.Range("J1:AA7").Clear()

# ---------------------------------------------------------------------------
# Sheet 3: "Lương" -- update a handful of labels / computed figures.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "Danh mục lương"
$ws3.Range("B2").Value = 17.5
$ws3.Range("B3").Value = 612500
$ws3.Range("B12").Value = 1875000
$ws3.Range("B29").Value = 2687500
$ws3.Range("B31").Value = 2687500
